$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Rupee (column B) values are stored as text (they look numeric but are
# shared-string cells), so force a Text number format before writing the new
# numeric-looking strings to keep Excel from auto-converting them to numbers.
$rng = $ws.Range("B2:B10")
$rng.NumberFormat = "@"

$ws.Range("B2").Value  = "732.91"
$ws.Range("B3").Value  = "3664.53"
$ws.Range("B4").Value  = "6596.15"
$ws.Range("B5").Value  = "8061.96"
$ws.Range("B6").Value  = "15391.01"
$ws.Range("B7").Value  = "2858.33"
$ws.Range("B8").Value  = "3517.94"
$ws.Range("B9").Value  = "4837.17"
$ws.Range("B10").Value = "2345.30"
